$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add Wins / Losses / Ties columns (AD, AE, AF) ---
# Copy formatting (bold font, border, center alignment) from the existing header cell AC1
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows 2-55: team record Wins=83, Losses=79, Ties=0 for every player row ---
$ws.Range("AD2:AD55").Value = 83
$ws.Range("AE2:AE55").Value = 79
$ws.Range("AF2:AF55").Value = 0
